# Update TPM-derived NATMI metrics (Fn1-Itga4 ligand-receptor pairs) for rows 2-10
# per commit "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.240107999999999
$ws.Range("H2").Value = 18.720324
$ws.Range("I2").Value = 0.01732230523539376
$ws.Range("J2").Value = 0.01732230523539376
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05968133333333333
$ws.Range("N2").Value = 0.179044
$ws.Range("O2").Value = 0.02602747651633847
$ws.Range("P2").Value = 0.02602747651633848
$ws.Range("Q2").Value = 0.3724179655839999
$ws.Range("R2").Value = 3.351761690255999
$ws.Range("S2").Value = 0.0004508558927230581
$ws.Range("T2").Value = 0.0004508558927230582

# Row 3
$ws.Range("G3").Value = 6.240107999999999
$ws.Range("H3").Value = 18.720324
$ws.Range("I3").Value = 0.01732230523539376
$ws.Range("J3").Value = 0.01732230523539376
$ws.Range("O3").Value = 0.144012433133819
$ws.Range("P3").Value = 0.144012433133819
$ws.Range("Q3").Value = 2.060623024012
$ws.Range("R3").Value = 18.545607216108
$ws.Range("S3").Value = 0.002494627324435747
$ws.Range("T3").Value = 0.002494627324435747

# Row 4
$ws.Range("G4").Value = 6.240107999999999
$ws.Range("H4").Value = 18.720324
$ws.Range("I4").Value = 0.01732230523539376
$ws.Range("J4").Value = 0.01732230523539376
$ws.Range("O4").Value = 0.8299600903498424
$ws.Range("P4").Value = 0.8299600903498425
$ws.Range("Q4").Value = 11.875605695772
$ws.Range("R4").Value = 106.880451261948
$ws.Range("S4").Value = 0.01437682201823495
$ws.Range("T4").Value = 0.01437682201823496

# Row 5
$ws.Range("I5").Value = 0.9592798330716089
$ws.Range("J5").Value = 0.9592798330716091
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05968133333333333
$ws.Range("N5").Value = 0.179044
$ws.Range("O5").Value = 0.02602747651633847
$ws.Range("P5").Value = 0.02602747651633848
$ws.Range("Q5").Value = 20.62387419015867
$ws.Range("R5").Value = 185.614867711428
$ws.Range("S5").Value = 0.02496763332786839
$ws.Range("T5").Value = 0.0249676333278684

# Row 6
$ws.Range("I6").Value = 0.9592798330716089
$ws.Range("J6").Value = 0.9592798330716091
$ws.Range("O6").Value = 0.144012433133819
$ws.Range("P6").Value = 0.144012433133819
$ws.Range("S6").Value = 0.1381482228168461
$ws.Range("T6").Value = 0.1381482228168462

# Row 7
$ws.Range("I7").Value = 0.9592798330716089
$ws.Range("J7").Value = 0.9592798330716091
$ws.Range("O7").Value = 0.8299600903498424
$ws.Range("P7").Value = 0.8299600903498425
$ws.Range("S7").Value = 0.7961639769268943
$ws.Range("T7").Value = 0.7961639769268946

# Row 8
$ws.Range("G8").Value = 8.428738666666666
$ws.Range("I8").Value = 0.02339786169299727
$ws.Range("J8").Value = 0.02339786169299728
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05968133333333333
$ws.Range("N8").Value = 0.179044
$ws.Range("O8").Value = 0.02602747651633847
$ws.Range("P8").Value = 0.02602747651633848
$ws.Range("Q8").Value = 0.5030383619448888
$ws.Range("R8").Value = 4.527345257504
$ws.Range("S8").Value = 0.000608987295747022
$ws.Range("T8").Value = 0.0006089872957470222

# Row 9
$ws.Range("G9").Value = 8.428738666666666
$ws.Range("I9").Value = 0.02339786169299727
$ws.Range("J9").Value = 0.02339786169299728
$ws.Range("O9").Value = 0.144012433133819
$ws.Range("P9").Value = 0.144012433133819
$ws.Range("Q9").Value = 2.783357749563555
$ws.Range("S9").Value = 0.003369582992537114
$ws.Range("T9").Value = 0.003369582992537115

# Row 10
$ws.Range("G10").Value = 8.428738666666666
$ws.Range("I10").Value = 0.02339786169299727
$ws.Range("J10").Value = 0.02339786169299728
$ws.Range("O10").Value = 0.8299600903498424
$ws.Range("P10").Value = 0.8299600903498425
$ws.Range("S10").Value = 0.01941929140471313
$ws.Range("T10").Value = 0.01941929140471314
